$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings keep exact formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.319.24'
$ws.Range("E2").Value = '  -0.39%  '

$ws.Range("D3").Value = '1.839.71'
$ws.Range("E3").Value = '  -0.91%  '

$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").Value = '239.99'
$ws.Range("E5").Value = '  -0.68%  '

$ws.Range("D6").Value = '0.6287'
$ws.Range("E6").Value = '  -0.85%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '0.07512'
$ws.Range("E8").Value = '  -1.09%  '

$ws.Range("D9").Value = '0.2899'
$ws.Range("E9").Value = '  -0.97%  '

$ws.Range("D10").Value = '24.78'
$ws.Range("E10").Value = '  +0.72%  '

$ws.Range("D11").Value = '0.07734'
$ws.Range("E11").Value = '  -0.34%  '

$ws.Range("D12").Value = '1.842.24'
$ws.Range("E12").Value = '  -0.73%  '

$ws.Range("D13").Value = '4.976'
$ws.Range("E13").Value = '  -1.18%  '

$ws.Range("D14").Value = '0.6779'
$ws.Range("E14").Value = '  -1.24%  '

$ws.Range("D15").Value = '0.00001021'
$ws.Range("E15").Value = '  -2.62%  '

$ws.Range("D17").Value = '6.240'
$ws.Range("E17").Value = '  +1.29%  '

$ws.Range("D18").Value = '29.323.62'
$ws.Range("E18").Value = '  -0.39%  '

$ws.Range("D19").Value = '228.89'
$ws.Range("E19").Value = '  -0.75%  '

$ws.Range("D20").Value = '12.31'
$ws.Range("E20").Value = '  -0.72%  '

$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").Value = '  -0.10%  '

$ws.Range("D22").Value = '7.417'
$ws.Range("E22").Value = '  -1.52%  '

$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").Value = '158.96'
$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("D25").Value = '8.463'
$ws.Range("E25").Value = '  -0.28%  '

$ws.Range("D26").Value = '0.1353'
$ws.Range("E26").Value = '  -3.69%  '

$ws.Range("D27").Value = '17.40'
$ws.Range("E27").Value = '  -1.95%  '

$ws.Range("B28").Value = 'Hedera'
$ws.Range("C28").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D28").Value = '0.06492'
$ws.Range("E28").Value = '  +13.69%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.447'
$ws.Range("E29").Value = '  +1.97%  '

$ws.Range("D30").Value = '1.485'
$ws.Range("E30").Value = '  +0.28%  '

$ws.Range("D31").Value = '4.065'
$ws.Range("E31").Value = '  -2.31%  '

$ws.Range("E32").Value = '  -0.12%  '

$ws.Range("D33").Value = '1.838'
$ws.Range("E33").Value = '  +0.29%  '

$ws.Range("D34").Value = '1.138'
$ws.Range("E34").Value = '  -1.84%  '

$ws.Range("D35").Value = '0.6958'
$ws.Range("E35").Value = '  -0.31%  '

$ws.Range("D36").Value = '2.570'
$ws.Range("E36").Value = '  -0.75%  '

$ws.Range("D37").Value = '0.01855'
$ws.Range("E37").Value = '  +1.51%  '

$ws.Range("D38").Value = '2.815'
$ws.Range("E38").Value = '  +1.55%  '

$ws.Range("D39").Value = '1.240.29'
$ws.Range("E39").Value = '  -0.85%  '

$ws.Range("D40").Value = '6.758'
$ws.Range("E40").Value = '  +3.51%  '

$ws.Range("D41").Value = '0.9301'
$ws.Range("E41").Value = '  +2.15%  '

$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '0.9993'
$ws.Range("E42").Value = '  -0.21%  '

$ws.Range("B43").Value = 'RocketPoolETH'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D43").Value = '1.993.92'
$ws.Range("E43").Value = '  -1.08%  '

$ws.Range("D44").Value = '100.83'
$ws.Range("E44").Value = '  -0.43%  '

$ws.Range("D45").Value = '65.55'
$ws.Range("E45").Value = '  -0.80%  '

$ws.Range("D46").Value = '0.00000000119'
$ws.Range("E46").Value = '  +3.61%  '

$ws.Range("D47").Value = '7.047'
$ws.Range("E47").Value = '  -1.72%  '

$ws.Range("D48").Value = '1.712'
$ws.Range("E48").Value = '  +1.69%  '

$ws.Range("D49").Value = '0.1153'
$ws.Range("E49").Value = '  -1.34%  '

$ws.Range("D50").Value = '9.004'
$ws.Range("E50").Value = '  -0.87%  '

$ws.Range("D51").Value = '0.3899'
$ws.Range("E51").Value = '  -1.89%  '

# Restore default styling on column D now that values are set
$ws.Range("D2:D51").Style = "Normal"

Write-Host "applied cryptos update"